# Prog-entera: Modifica datos del problema 11.
# Updates the parameter values on the "Otros-parametros" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Otros-parametros")

$ws.Range("B2").Value = 210
$ws.Range("B3").Value = 20
$ws.Range("B4").Value = 3200000
$ws.Range("B5").Value = 6400000
$ws.Range("B6").Value = 9600000
